$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "2026-01-31 10:08"
$ws.Range("B12").Value = 23
$ws.Range("C12").Value = 5
